$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Debug_Messages" setting row
$ws.Range("A3").Value = "Debug_Messages"
$ws.Range("B3").Value = 0

# Move the active selection to F6 (matches the saved view state in the diff)
$null = $ws.Range("F6").Select()

# Update the workbook window position (best-effort view state change;
# this is a cosmetic Excel-window placement value)
$excel.ActiveWindow.Left = 3040
$excel.ActiveWindow.Top = 3040
